# "swati addded new classes"
#
# On the "ListManager_BlackList" sheet, cell C2 held the mailto hyperlink
# "swati@gmail.com". Swati edited it to a new address, "swati45@gmail.com",
# while keeping the hyperlink (its original address text is preserved as the
# link's "display" text), and then left her selection on that sheet
# (cell C11), which is now the active tab in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListManager_BlackList")

$cell = $ws.Range("C2")

# Recreate the existing hyperlink on C2 so the original address text can be
# preserved as the link's display text once the cell text itself changes.
$cell.Hyperlinks.Delete()
$ws.Hyperlinks.Add($cell, "mailto:swati@gmail.com", $null, $null, "swati@gmail.com")

# New email address for this row (leading "'" keeps it as plain text).
$cell.Value = "'swati45@gmail.com"

# Swati ends up on the BlackList sheet with C11 selected.
$ws.Activate()
$ws.Range("C11").Select()
